$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-07-03 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-04 Friday", 2)

# Update the multiplication problems in the table, cell by cell so that the
# two duplicate "74x18=" entries each map to their own distinct replacement.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "71×26="
$tbl.Cell(1, 2).Range.Text = "88×11="
$tbl.Cell(1, 3).Range.Text = "38×94="
$tbl.Cell(1, 4).Range.Text = "35×63="
$tbl.Cell(1, 5).Range.Text = "85×47="

$tbl.Cell(5, 1).Range.Text = "91×17="
$tbl.Cell(5, 2).Range.Text = "27×45="
$tbl.Cell(5, 3).Range.Text = "94×72="
$tbl.Cell(5, 4).Range.Text = "13×87="
$tbl.Cell(5, 5).Range.Text = "62×73="

$tbl.Cell(10, 1).Range.Text = "85×63="
$tbl.Cell(10, 2).Range.Text = "50×94="
$tbl.Cell(10, 3).Range.Text = "80×71="
$tbl.Cell(10, 4).Range.Text = "50×36="
$tbl.Cell(10, 5).Range.Text = "37×76="

$tbl.Cell(15, 1).Range.Text = "16×67="
$tbl.Cell(15, 2).Range.Text = "20×14="
$tbl.Cell(15, 3).Range.Text = "43×22="
$tbl.Cell(15, 4).Range.Text = "69×83="
$tbl.Cell(15, 5).Range.Text = "72×52="

$tbl.Cell(20, 1).Range.Text = "63×11="
$tbl.Cell(20, 2).Range.Text = "97×82="
$tbl.Cell(20, 3).Range.Text = "67×20="
$tbl.Cell(20, 4).Range.Text = "95×77="
$tbl.Cell(20, 5).Range.Text = "14×84="
